$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text prefix ensures values like "533.02" or "6.40" are stored as
# literal text (matching the source inlineStr cells) rather than being
# auto-converted to numbers by Excel's type inference. Resetting the style
# back to "Normal" afterwards drops the quote-prefix cell style so the
# cell keeps its original (unstyled) formatting.

$ws.Range("D2").Value = "'59.600.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.04%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.649.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.12%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'533.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.84%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'147.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.61%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.20%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.570"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.36%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -3.84%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +1.39%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.59%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.62%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.119.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.22%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'59.547.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.88%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -1.43%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +0.43%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.669.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.60%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'344.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.96%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.06%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +2.80%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.11%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.01%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'66.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.78%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +1.07%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.169"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.06%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.769.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.26%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.35%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.32%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0₃0803"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.16%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.05%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'6.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.29%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +1.92%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'19.06"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.45%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'150.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.98%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'4.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.49%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.55%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.865"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -4.07%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.859"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.99%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.04%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.32%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +0.87%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.996"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.34%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.0981"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.75%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.603"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.43%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'271.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.59%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'19.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.04%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +1.80%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0536"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.68%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.041.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.42%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'4.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.21%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0230"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.86%  "
$ws.Range("E51").Style = "Normal"
